# Apply updated betting-odds values for Jogos_da_Semana_FlashScore_2025-06-04.xlsx
# (rows/columns identified per the source diff; only cells whose value actually
# changed are touched, addressed via Cells.Item(row, col) for AA.. columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.6  # G2
$ws.Cells.Item(2, 9).Value = 6.5  # I2
$ws.Cells.Item(2, 31).Value = 12  # AE2

# Row 3
$ws.Cells.Item(3, 7).Value = 2.27  # G3
$ws.Cells.Item(3, 8).Value = 3.65  # H3
$ws.Cells.Item(3, 9).Value = 2.8  # I3
$ws.Cells.Item(3, 11).Value = 9.5  # K3
$ws.Cells.Item(3, 12).Value = 1.19  # L3
$ws.Cells.Item(3, 13).Value = 4.35  # M3
$ws.Cells.Item(3, 14).Value = 1.57  # N3
$ws.Cells.Item(3, 15).Value = 2.3  # O3
$ws.Cells.Item(3, 16).Value = 1.31  # P3
$ws.Cells.Item(3, 17).Value = 3.3  # Q3
$ws.Cells.Item(3, 18).Value = 1.5  # R3
$ws.Cells.Item(3, 19).Value = 2.42  # S3
$ws.Cells.Item(3, 20).Value = 10.75  # T3
$ws.Cells.Item(3, 21).Value = 15  # U3
$ws.Cells.Item(3, 22).Value = 9.75  # V3
$ws.Cells.Item(3, 24).Value = 17  # X3
$ws.Cells.Item(3, 25).Value = 22  # Y3
$ws.Cells.Item(3, 26).Value = 9.5  # Z3
$ws.Cells.Item(3, 27).Value = 7.8  # AA3
$ws.Cells.Item(3, 29).Value = 40  # AC3
$ws.Cells.Item(3, 30).Value = 250  # AD3
$ws.Cells.Item(3, 31).Value = 11.5  # AE3
$ws.Cells.Item(3, 32).Value = 18  # AF3
$ws.Cells.Item(3, 33).Value = 11  # AG3
$ws.Cells.Item(3, 35).Value = 22  # AI3
$ws.Cells.Item(3, 36).Value = 25  # AJ3

# Row 4
$ws.Cells.Item(4, 7).Value = 2.55  # G4
$ws.Cells.Item(4, 9).Value = 3  # I4
$ws.Cells.Item(4, 10).Value = 1.1  # J4
$ws.Cells.Item(4, 11).Value = 7  # K4
$ws.Cells.Item(4, 21).Value = 11  # U4
$ws.Cells.Item(4, 22).Value = 11  # V4
$ws.Cells.Item(4, 23).Value = 26  # W4
$ws.Cells.Item(4, 28).Value = 17  # AB4
$ws.Cells.Item(4, 32).Value = 13  # AF4
$ws.Cells.Item(4, 34).Value = 29  # AH4

# Row 5
$ws.Cells.Item(5, 7).Value = 2.05  # G5
$ws.Cells.Item(5, 9).Value = 4  # I5
$ws.Cells.Item(5, 18).Value = 2.1  # R5
$ws.Cells.Item(5, 19).Value = 1.67  # S5
$ws.Cells.Item(5, 21).Value = 8.5  # U5
$ws.Cells.Item(5, 23).Value = 17  # W5
$ws.Cells.Item(5, 31).Value = 9.5  # AE5
$ws.Cells.Item(5, 32).Value = 19  # AF5
$ws.Cells.Item(5, 33).Value = 15  # AG5
$ws.Cells.Item(5, 35).Value = 41  # AI5

# Row 12
$ws.Cells.Item(12, 7).Value = 2.2  # G12
$ws.Cells.Item(12, 8).Value = 3  # H12
$ws.Cells.Item(12, 9).Value = 3.35  # I12
$ws.Cells.Item(12, 12).Value = 1.5  # L12
$ws.Cells.Item(12, 13).Value = 2.25  # M12
$ws.Cells.Item(12, 14).Value = 2.45  # N12
$ws.Cells.Item(12, 15).Value = 1.42  # O12
$ws.Cells.Item(12, 16).Value = 1.53  # P12
$ws.Cells.Item(12, 17).Value = 2.18  # Q12
$ws.Cells.Item(12, 18).Value = 2.1  # R12
$ws.Cells.Item(12, 19).Value = 1.57  # S12
$ws.Cells.Item(12, 20).Value = 5.6  # T12
$ws.Cells.Item(12, 21).Value = 9  # U12
$ws.Cells.Item(12, 23).Value = 21  # W12
$ws.Cells.Item(12, 24).Value = 23  # X12
$ws.Cells.Item(12, 26).Value = 6.3  # Z12
$ws.Cells.Item(12, 27).Value = 6  # AA12
$ws.Cells.Item(12, 28).Value = 19.5  # AB12
$ws.Cells.Item(12, 29).Value = 150  # AC12
$ws.Cells.Item(12, 31).Value = 7.3  # AE12
$ws.Cells.Item(12, 32).Value = 15.5  # AF12
$ws.Cells.Item(12, 33).Value = 13  # AG12
$ws.Cells.Item(12, 35).Value = 40  # AI12
$ws.Cells.Item(12, 36).Value = 60  # AJ12

# Row 14
$ws.Cells.Item(14, 7).Value = 2.7  # G14
$ws.Cells.Item(14, 8).Value = 2.75  # H14
$ws.Cells.Item(14, 9).Value = 2.82  # I14
$ws.Cells.Item(14, 12).Value = 1.52  # L14
$ws.Cells.Item(14, 13).Value = 2.22  # M14
$ws.Cells.Item(14, 14).Value = 2.47  # N14
$ws.Cells.Item(14, 15).Value = 1.42  # O14
$ws.Cells.Item(14, 16).Value = 1.55  # P14
$ws.Cells.Item(14, 17).Value = 2.15  # Q14
$ws.Cells.Item(14, 18).Value = 2.02  # R14
$ws.Cells.Item(14, 19).Value = 1.62  # S14
$ws.Cells.Item(14, 20).Value = 6.6  # T14
$ws.Cells.Item(14, 21).Value = 12.5  # U14
$ws.Cells.Item(14, 22).Value = 10.5  # V14
$ws.Cells.Item(14, 23).Value = 32  # W14
$ws.Cells.Item(14, 24).Value = 28  # X14
$ws.Cells.Item(14, 25).Value = 45  # Y14
$ws.Cells.Item(14, 26).Value = 5.9  # Z14
$ws.Cells.Item(14, 27).Value = 5.5  # AA14
$ws.Cells.Item(14, 28).Value = 17.5  # AB14
$ws.Cells.Item(14, 29).Value = 110  # AC14
$ws.Cells.Item(14, 31).Value = 6.3  # AE14
$ws.Cells.Item(14, 32).Value = 12.5  # AF14
$ws.Cells.Item(14, 33).Value = 11.25  # AG14
$ws.Cells.Item(14, 34).Value = 35  # AH14
$ws.Cells.Item(14, 35).Value = 32  # AI14
$ws.Cells.Item(14, 36).Value = 50  # AJ14

# Row 16
$ws.Cells.Item(16, 7).Value = 2.72  # G16
$ws.Cells.Item(16, 8).Value = 3.55  # H16
$ws.Cells.Item(16, 9).Value = 2.32  # I16
$ws.Cells.Item(16, 10).Value = 1.04  # J16
$ws.Cells.Item(16, 11).Value = 8.5  # K16
$ws.Cells.Item(16, 12).Value = 1.22  # L16
$ws.Cells.Item(16, 13).Value = 3.9  # M16
$ws.Cells.Item(16, 14).Value = 1.65  # N16
$ws.Cells.Item(16, 15).Value = 2.1  # O16
$ws.Cells.Item(16, 16).Value = 1.34  # P16
$ws.Cells.Item(16, 17).Value = 3  # Q16
$ws.Cells.Item(16, 18).Value = 1.55  # R16
$ws.Cells.Item(16, 19).Value = 2.3  # S16
$ws.Cells.Item(16, 20).Value = 11  # T16
$ws.Cells.Item(16, 21).Value = 15.5  # U16
$ws.Cells.Item(16, 22).Value = 10  # V16
$ws.Cells.Item(16, 23).Value = 32  # W16
$ws.Cells.Item(16, 24).Value = 20  # X16
$ws.Cells.Item(16, 25).Value = 25  # Y16
$ws.Cells.Item(16, 26).Value = 8.5  # Z16
$ws.Cells.Item(16, 27).Value = 7  # AA16
$ws.Cells.Item(16, 28).Value = 12  # AB16
$ws.Cells.Item(16, 29).Value = 45  # AC16
$ws.Cells.Item(16, 30).Value = 250  # AD16
$ws.Cells.Item(16, 31).Value = 10.25  # AE16
$ws.Cells.Item(16, 32).Value = 13  # AF16
$ws.Cells.Item(16, 33).Value = 9  # AG16
$ws.Cells.Item(16, 34).Value = 25  # AH16
$ws.Cells.Item(16, 35).Value = 17  # AI16
$ws.Cells.Item(16, 36).Value = 22  # AJ16
